{"js": "// Update the division-problem answer table: each data row (every 4th row,\n// starting at row 0) holds 5 \"A\u00f7B=C, D\" strings that need new values.\nconst newValuesByRow = {\n  0: [\"911\u00f77=130, 1\", \"438\u00f73=146, 0\", \"605\u00f79=67, 2\", \"858\u00f72=429, 0\", \"750\u00f77=107, 1\"],\n  4: [\"878\u00f76=146, 2\", \"193\u00f74=48, 1\", \"189\u00f78=23, 5\", \"503\u00f78=62, 7\", \"862\u00f79=95, 7\"],\n  8: [\"665\u00f76=110, 5\", \"686\u00f79=76, 2\", \"996\u00f79=110, 6\", \"654\u00f77=93, 3\", \"965\u00f77=137, 6\"],\n  12: [\"940\u00f73=313, 1\", \"292\u00f75=58, 2\", \"942\u00f76=157, 0\", \"758\u00f78=94, 6\", \"777\u00f75=155, 2\"],\n  16: [\"820\u00f75=164, 0\", \"136\u00f74=34, 0\", \"459\u00f75=91, 4\", \"582\u00f72=291, 0\", \"979\u00f79=108, 7\"],\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const rowIndexStr of Object.keys(newValuesByRow)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const rowValues = newValuesByRow[rowIndexStr];\n  for (let colIndex = 0; colIndex < rowValues.length; colIndex++) {\n    const cell = table.getCell(rowIndex, colIndex);\n    cell.value = rowValues[colIndex];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem answer table: each data row (every 4th row,\n# starting at row 1) holds 5 \"A\u00f7B=C, D\" strings that need new values.\n$doc = $word.ActiveDocument\n$table = $doc.Tables.Item(1)\n\n$rowsMap = @{\n    1  = @(\"911\u00f77=130, 1\", \"438\u00f73=146, 0\", \"605\u00f79=67, 2\", \"858\u00f72=429, 0\", \"750\u00f77=107, 1\")\n    5  = @(\"878\u00f76=146, 2\", \"193\u00f74=48, 1\", \"189\u00f78=23, 5\", \"503\u00f78=62, 7\", \"862\u00f79=95, 7\")\n    9  = @(\"665\u00f76=110, 5\", \"686\u00f79=76, 2\", \"996\u00f79=110, 6\", \"654\u00f77=93, 3\", \"965\u00f77=137, 6\")\n    13 = @(\"940\u00f73=313, 1\", \"292\u00f75=58, 2\", \"942\u00f76=157, 0\", \"758\u00f78=94, 6\", \"777\u00f75=155, 2\")\n    17 = @(\"820\u00f75=164, 0\", \"136\u00f74=34, 0\", \"459\u00f75=91, 4\", \"582\u00f72=291, 0\", \"979\u00f79=108, 7\")\n}\n\nforeach ($rowIndex in $rowsMap.Keys) {\n    $vals = $rowsMap[$rowIndex]\n    for ($c = 1; $c -le $vals.Length; $c++) {\n        $cell = $table.Cell($rowIndex, $c)\n        $cell.Range.Text = $vals[$c - 1]\n    }\n}\n"}
